$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8606349999999999
$ws.Range("H2").Value = 2.581905
$ws.Range("I2").Value = 0.0262626340301864
$ws.Range("J2").Value = 0.0262626340301864
$ws.Range("M2").Value = 135.0916853333333
$ws.Range("N2").Value = 405.2750559999999
$ws.Range("O2").Value = 0.7123704212620513
$ws.Range("P2").Value = 0.7123704212620514
$ws.Range("Q2").Value = 116.2646326068533
$ws.Range("R2").Value = 1046.38169346168
$ws.Range("S2").Value = 0.01870872366753497
$ws.Range("T2").Value = 0.01870872366753497

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8606349999999999
$ws.Range("H3").Value = 2.581905
$ws.Range("I3").Value = 0.0262626340301864
$ws.Range("J3").Value = 0.0262626340301864
$ws.Range("O3").Value = 0.2125756143240238
$ws.Range("P3").Value = 0.2125756143240238
$ws.Range("Q3").Value = 34.69406500170666
$ws.Range("R3").Value = 312.24658501536
$ws.Range("S3").Value = 0.005582795562733886
$ws.Range("T3").Value = 0.005582795562733886

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8606349999999999
$ws.Range("H4").Value = 2.581905
$ws.Range("I4").Value = 0.0262626340301864
$ws.Range("J4").Value = 0.0262626340301864
$ws.Range("M4").Value = 14.23299766666667
$ws.Range("N4").Value = 42.698993
$ws.Range("O4").Value = 0.07505396441392481
$ws.Range("P4").Value = 0.07505396441392483
$ws.Range("Q4").Value = 12.24941594685167
$ws.Range("R4").Value = 110.244743521665
$ws.Range("S4").Value = 0.001971114799917541
$ws.Range("T4").Value = 0.001971114799917541

# Row 5
$ws.Range("H5").Value = 58.40949000000001
$ws.Range("I5").Value = 0.5941299388474139
$ws.Range("J5").Value = 0.5941299388474139
$ws.Range("M5").Value = 135.0916853333333
$ws.Range("N5").Value = 405.2750559999999
$ws.Range("O5").Value = 0.7123704212620513
$ws.Range("P5").Value = 0.7123704212620514
$ws.Range("Q5").Value = 2630.212147853493
$ws.Range("R5").Value = 23671.90933068144
$ws.Range("S5").Value = 0.423240594821129
$ws.Range("T5").Value = 0.4232405948211291

# Row 6
$ws.Range("H6").Value = 58.40949000000001
$ws.Range("I6").Value = 0.5941299388474139
$ws.Range("J6").Value = 0.5941299388474139
$ws.Range("O6").Value = 0.2125756143240238
$ws.Range("P6").Value = 0.2125756143240238
$ws.Range("Q6").Value = 784.8711098109867
$ws.Range("R6").Value = 7063.839988298881
$ws.Range("S6").Value = 0.1262975367387837
$ws.Range("T6").Value = 0.1262975367387837

# Row 7
$ws.Range("H7").Value = 58.40949000000001
$ws.Range("I7").Value = 0.5941299388474139
$ws.Range("J7").Value = 0.5941299388474139
$ws.Range("M7").Value = 14.23299766666667
$ws.Range("N7").Value = 42.698993
$ws.Range("O7").Value = 0.07505396441392481
$ws.Range("P7").Value = 0.07505396441392483
$ws.Range("Q7").Value = 277.1140449603967
$ws.Range("R7").Value = 2494.02640464357
$ws.Range("S7").Value = 0.04459180728750113
$ws.Range("T7").Value = 0.04459180728750114

# Row 8
$ws.Range("G8").Value = 12.439858
$ws.Range("H8").Value = 37.319574
$ws.Range("I8").Value = 0.3796074271223998
$ws.Range("J8").Value = 0.3796074271223997
$ws.Range("M8").Value = 135.0916853333333
$ws.Range("N8").Value = 405.2750559999999
$ws.Range("O8").Value = 0.7123704212620513
$ws.Range("P8").Value = 0.7123704212620514
$ws.Range("Q8").Value = 1680.521382527349
$ws.Range("R8").Value = 15124.69244274614
$ws.Range("S8").Value = 0.2704211027733874
$ws.Range("T8").Value = 0.2704211027733874

# Row 9
$ws.Range("G9").Value = 12.439858
$ws.Range("H9").Value = 37.319574
$ws.Range("I9").Value = 0.3796074271223998
$ws.Range("J9").Value = 0.3796074271223997
$ws.Range("O9").Value = 0.2125756143240238
$ws.Range("P9").Value = 0.2125756143240238
$ws.Range("Q9").Value = 501.4776787650986
$ws.Range("R9").Value = 4513.299108885887
$ws.Range("S9").Value = 0.08069528202250621
$ws.Range("T9").Value = 0.08069528202250621

# Row 10
$ws.Range("G10").Value = 12.439858
$ws.Range("H10").Value = 37.319574
$ws.Range("I10").Value = 0.3796074271223998
$ws.Range("J10").Value = 0.3796074271223997
$ws.Range("M10").Value = 14.23299766666667
$ws.Range("N10").Value = 42.698993
$ws.Range("O10").Value = 0.07505396441392481
$ws.Range("P10").Value = 0.07505396441392483
$ws.Range("Q10").Value = 177.0564698876647
$ws.Range("R10").Value = 1593.508228988982
$ws.Range("S10").Value = 0.02849104232650615
$ws.Range("T10").Value = 0.02849104232650615

